# Apply the "error-various-content" fixture edits described in the commit.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("registrants")
$ws2 = $wb.Worksheets.Item("services")

# --- Sheet "registrants" ---------------------------------------------------
# Row 2 (Adrien): remove the last name entirely and flag the missing data.
$ws1.Range("B2").ClearContents()
$ws1.Range("E2").Value = "Where is the lastname ?!"

# Row 4 (Sophie Roekhaut): her interface language is bogus ("BOUP"), flag it.
$ws1.Range("D4").Value = "BOUP"
$ws1.Range("E4").Value = "BOUP is not even a language"

# Row 6 (..Agnerey): interface language corrected from PT_BR to IT_IT.
$ws1.Range("D6").Value = "IT_IT"

# --- Sheet "services" -------------------------------------------------------
# Row 3: now belongs to Sophie (sroekhaut), SRS course in DE_DE, 45 days.
$ws2.Range("A3").Value = "sroekhaut@altissia.org"
$ws2.Range("C3").Value = "DE_DE"
$ws2.Range("D3").Value = 45

# Row 4: Sophie also has a COURSE in NL_NL for 60 days.
$ws2.Range("B4").Value = "COURSE"
$ws2.Range("C4").Value = "NL_NL"
$ws2.Range("D4").Value = 60

# Row 5: now belongs to Renaud (rlaloux), COURSE in EN_GB, 10 days.
$ws2.Range("A5").Value = "rlaloux@altissia.org"
$ws2.Range("C5").Value = "EN_GB"
$ws2.Range("D5").Value = 10

# Row 6: Renaud also has an ASSESSMENT in EN_GB, duration cleared.
$ws2.Range("B6").Value = "ASSESSMENT"
$ws2.Range("C6").Value = "EN_GB"
$ws2.Range("D6").ClearContents()

# Row 7: now belongs to mbacoup, COURSE in FR_FR, 50 days, missing-registrant note.
$ws2.Range("A7").Value = "mbacoup@beable.com"
$ws2.Range("B7").Value = "COURSE"
$ws2.Range("C7").Value = "FR_FR"
$ws2.Range("D7").Value = 50
$ws2.Range("E7").Value = "mbacoup is absent from the registrants sheet"

# Row 8: now belongs to gagneray, NEWS in DE_DE, -1 duration, "should at least be zero" note.
$ws2.Range("A8").Value = "gagneray@beable.com"
$ws2.Range("B8").Value = "NEWS"
$ws2.Range("C8").Value = "DE_DE"
$ws2.Range("D8").Value = -1
$ws2.Range("E8").Value = "should at least be zero"

# --- View state: activate "registrants" and select E3, deselect "services" tab ---
$ws2.Range("B4").Select()
$ws1.Activate()
$ws1.Range("E3").Select()
